# remise du projet final
#
# PowerPoint "Designer" style metadata stamp: every shape, on every
# slide, gets a custom "NUM" tag recording its 1-based position within
# the slide's shape tree (this is what PowerPoint itself writes to
# ppt/tags/tagN.xml + <p:custDataLst><p:tags r:id="..."/></p:custDataLst>
# inside each shape's <p:nvPr> when the Designer pane has touched the
# slide).

$p = $ppt.ActivePresentation

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        $shape.Tags.Add("NUM", "$shapeIdx")
    }
}
